$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.742.16'
$ws.Range('E2').Value = '  -1.08%  '
$ws.Range('D3').Value = '1.623.10'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '214.79'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5069'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.24%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.002'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2559'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.09%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.24'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.27%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07773'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.245'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.26%  '
$ws.Range('D13').Value = '1.623.92'
$ws.Range('E13').Value = '  -1.12%  '
$ws.Range('D14').Value = '1.847.07'
$ws.Range('E14').Value = '  -1.08%  '
$ws.Range('E15').Value = '  +1.09%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '63.54'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.71%  '
$ws.Range('D17').Value = '0.0₅7522'
$ws.Range('E17').Value = '  -3.04%  '
$ws.Range('D18').Value = '25.765.80'
$ws.Range('E18').Value = '  -1.02%  '
$ws.Range('E19').Value = '  -0.33%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '193.45'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.61%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.387'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.34%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.752'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.40%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.954'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.58%  '
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.865'
$ws.Range('D25').ClearFormats()
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '140.68'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.10%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1240'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.92%  '
$ws.Range('E28').Value = '  -2.12%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.43'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.81%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.234'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.42%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.04859'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.07%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.309'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.64%  '
$ws.Range('E33').Value = '  -1.53%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.545'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.363'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.56%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.8917'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.93%  '
$ws.Range('D37').Value = '1.127.74'
$ws.Range('E37').Value = '  +1.35%  '
$ws.Range('E38').Value = '  -2.01%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.5490'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01554'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.16%  '
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.567'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.54%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.7944'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.11%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '97.05'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.56%  '
$ws.Range('D45').Value = '1.770.46'
$ws.Range('E45').Value = '  -0.38%  '
$ws.Range('D46').Value = '0.0₈114'
$ws.Range('E46').Value = '  -8.26%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4419'
$ws.Range('D47').ClearFormats()
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '54.64'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.96%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.05126'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.90%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.594'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.96%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.9987'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.95%  '
